# The workbook originally has a single sheet "Sheet1" holding a small
# hierarchical issue-tracker table (outline/grouped rows under an "Key"
# column). This edit:
#   1. Renames the existing sheet "Sheet1" -> "Sheet2" and empties it.
#   2. Adds a new sheet "to_import" right after it, carrying the table
#      that used to live on the original sheet (same values + row
#      outline/grouping), with the header "Key" renamed to "Id".

$wb = $excel.ActiveWorkbook

# --- Step 1: grab the existing (only) sheet and rename it -----------------
$oldSheet = $wb.Worksheets.Item(1)
$oldSheet.Name = "Sheet2"

# --- Step 2: add the new sheet right after the renamed one -----------------
$newSheet = $wb.Worksheets.Add($null, $oldSheet)
$newSheet.Name = "to_import"

# --- Step 3: write the table (former Sheet1 contents) onto to_import ------
# header row
$newSheet.Range("A1").Value = "Id"
$newSheet.Range("B1").Value = "Sub-item count"
$newSheet.Range("C1").Value = "Description"

# data rows: (row, key value, description, outline level)
$rows = @(
    @(2,  1,  "issue 1",  0),
    @(3,  2,  "issue 2",  0),
    @(4,  3,  "issue 3",  0),
    @(5,  4,  "issue 4",  1),
    @(6,  5,  "issue 5",  2),
    @(7,  6,  "issue 6",  2),
    @(8,  7,  "issue 7",  1),
    @(9,  8,  "issue 8",  0),
    @(10, 9,  "issue 9",  0),
    @(11, 10, "issue 10", 0)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $keyVal = $r[1]
    $desc = $r[2]
    $level = $r[3]

    $newSheet.Range("A$rowNum").Value = $keyVal
    $newSheet.Range("C$rowNum").Value = $desc

    if ($level -gt 0) {
        $newSheet.Rows($rowNum).OutlineLevel = $level
    }
}

# --- Step 4: clear the old sheet (now "Sheet2") so it is blank ------------
$oldSheet.Cells.Clear()
foreach ($r in $rows) {
    $rowNum = $r[0]
    $level = $r[3]
    if ($level -gt 0) {
        $oldSheet.Rows($rowNum).OutlineLevel = 0
    }
}

# Put the selection back on the first sheet, matching a freshly-opened book.
$wb.Worksheets.Item(1).Activate()
